$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-07 02:17:41'
$ws.Range('N2').Value = '-1.5 °C 1:50 TU'
$ws.Range('O2').Value = '-1.1 °C'
$ws.Range('E3').Value = '2026-02-07 02:17:44'
$ws.Range('N3').Value = '-5.9 °C 1:42 TU'
$ws.Range('O3').Value = '-4.7 °C'
$ws.Range('E4').Value = '2026-02-07 02:17:46'
$ws.Range('J4').Value = '1000.7 hPa'
$ws.Range('L4').Value = '35.3 km/h - 265º 1:52 TU'
$ws.Range('E5').Value = '2026-02-07 02:17:49'
$ws.Range('H5').Value = '''70%'
$ws.Range('O5').Value = '9.5 °C'
$ws.Range('E6').Value = '2026-02-07 02:17:51'
$ws.Range('H6').Value = '''57%'
$ws.Range('J6').Value = '1002.0 hPa'
$ws.Range('N6').Value = '11.4 °C 1:56 TU'
$ws.Range('O6').Value = '12.3 °C'
$ws.Range('E7').Value = '2026-02-07 02:17:54'
$ws.Range('H7').Value = '''72%'
$ws.Range('J7').Value = '1001.9 hPa'
$ws.Range('N7').Value = '7.7 °C 1:40 TU'
$ws.Range('O7').Value = '8.4 °C'
$ws.Range('E8').Value = '2026-02-07 02:17:56'
$ws.Range('L8').Value = '5.4 km/h - 75º 1:51 TU'
$ws.Range('N8').Value = '4.7 °C 1:56 TU'
$ws.Range('O8').Value = '5.2 °C'
$ws.Range('E9').Value = '2026-02-07 02:17:58'
$ws.Range('N9').Value = '2.4 °C 1:57 TU'
$ws.Range('O9').Value = '2.9 °C'
$ws.Range('E10').Value = '2026-02-07 02:18:01'
$ws.Range('M10').Value = '7.8 °C 1:45 TU'
$ws.Range('O10').Value = '7.1 °C'
$ws.Range('E11').Value = '2026-02-07 02:18:03'
$ws.Range('H11').Value = '''96%'
$ws.Range('J11').Value = '1005.1 hPa'
$ws.Range('E12').Value = '2026-02-07 02:18:05'
$ws.Range('L12').Value = '35.6 km/h - 283º 1:45 TU'
$ws.Range('O12').Value = '10.6 °C'
$ws.Range('E13').Value = '2026-02-07 02:18:08'
$ws.Range('O13').Value = '7.1 °C'
$ws.Range('E14').Value = '2026-02-07 02:18:10'
$ws.Range('H14').Value = '''88%'
$ws.Range('L14').Value = '18.4 km/h - 218º 1:40 TU'
$ws.Range('O14').Value = '-5.7 °C'
$ws.Range('E15').Value = '2026-02-07 02:18:13'
$ws.Range('H15').Value = '''77%'
$ws.Range('E16').Value = '2026-02-07 02:18:15'
$ws.Range('H16').Value = '''85%'
$ws.Range('N16').Value = '3.3 °C 1:54 TU'
$ws.Range('O16').Value = '4.0 °C'
$ws.Range('E17').Value = '2026-02-07 02:18:18'
$ws.Range('H17').Value = '''96%'
$ws.Range('O17').Value = '3.7 °C'
$ws.Range('E18').Value = '2026-02-07 02:18:20'
$ws.Range('N18').Value = '-6.5 °C 1:59 TU'
$ws.Range('O18').Value = '-6.2 °C'
$ws.Range('E19').Value = '2026-02-07 02:18:23'
$ws.Range('J19').Value = '1005.1 hPa'
$ws.Range('N19').Value = '4.8 °C 1:36 TU'
$ws.Range('O19').Value = '5.3 °C'
$ws.Range('E20').Value = '2026-02-07 02:18:25'
$ws.Range('N20').Value = '-4.8 °C 1:42 TU'
$ws.Range('O20').Value = '-4.2 °C'
$ws.Range('E21').Value = '2026-02-07 02:18:28'
$ws.Range('H21').Value = '''61%'
$ws.Range('J21').Value = '1000.6 hPa'
$ws.Range('N21').Value = '7.4 °C 1:50 TU'
$ws.Range('O21').Value = '9.8 °C'
$ws.Range('E22').Value = '2026-02-07 02:18:30'
$ws.Range('O22').Value = '6.2 °C'
$ws.Range('E23').Value = '2026-02-07 02:18:33'
$ws.Range('N23').Value = '7.6 °C 1:46 TU'
$ws.Range('E24').Value = '2026-02-07 02:18:35'
$ws.Range('N24').Value = '10.0 °C 1:36 TU'
$ws.Range('O24').Value = '10.6 °C'
$ws.Range('E25').Value = '2026-02-07 02:18:37'
$ws.Range('H25').Value = '''96%'
$ws.Range('J25').Value = '1004.7 hPa'
$ws.Range('N25').Value = '0.3 °C 1:30 TU'
$ws.Range('O25').Value = '1.0 °C'
$ws.Range('E26').Value = '2026-02-07 02:18:40'
$ws.Range('G26').Value = '111 cm'
$ws.Range('H26').Value = '''78%'
$ws.Range('L26').Value = '31.7 km/h - 39º 1:48 TU'
$ws.Range('M26').Value = '-0.6 °C 1:58 TU'
$ws.Range('O26').Value = '-1.3 °C'
$ws.Range('E27').Value = '2026-02-07 02:18:42'
$ws.Range('H27').Value = '''97%'
$ws.Range('N27').Value = '7.7 °C 1:36 TU'
$ws.Range('O27').Value = '8.3 °C'
$ws.Range('E28').Value = '2026-02-07 02:18:45'
$ws.Range('J28').Value = '1002.8 hPa'
$ws.Range('N28').Value = '3.6 °C 1:38 TU'
$ws.Range('O28').Value = '4.4 °C'
$ws.Range('E29').Value = '2026-02-07 02:18:47'
$ws.Range('H29').Value = '''54%'
$ws.Range('N29').Value = '10.5 °C 1:59 TU'
$ws.Range('O29').Value = '12.2 °C'
$ws.Range('E30').Value = '2026-02-07 02:18:50'
$ws.Range('H30').Value = '''80%'
$ws.Range('I30').Value = '0.4 mm'
$ws.Range('N30').Value = '-5.1 °C 1:52 TU'
$ws.Range('O30').Value = '-4.3 °C'
$ws.Range('E31').Value = '2026-02-07 02:18:52'
$ws.Range('N31').Value = '3.7 °C 1:52 TU'
$ws.Range('O31').Value = '3.9 °C'
$ws.Range('E32').Value = '2026-02-07 02:18:54'
$ws.Range('H32').Value = '''63%'
$ws.Range('J32').Value = '1003.4 hPa'
$ws.Range('L32').Value = '31.0 km/h - 277º 1:35 TU'
$ws.Range('O32').Value = '11.6 °C'
$ws.Range('E33').Value = '2026-02-07 02:18:57'
$ws.Range('H33').Value = '''92%'
$ws.Range('M33').Value = '9.5 °C 1:56 TU'
$ws.Range('O33').Value = '7.5 °C'
$ws.Range('E34').Value = '2026-02-07 02:18:59'
$ws.Range('H34').Value = '''73%'
$ws.Range('L34').Value = '19.1 km/h - 270º 1:56 TU'
$ws.Range('O34').Value = '7.3 °C'
$ws.Range('E35').Value = '2026-02-07 02:19:02'
$ws.Range('N35').Value = '-4.1 °C 1:30 TU'
$ws.Range('E36').Value = '2026-02-07 02:19:04'
$ws.Range('J36').Value = '1005.6 hPa'
$ws.Range('L36').Value = '10.1 km/h - 52º 1:53 TU'

Write-Host "Applied all changes"
